$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set target cells to text format to preserve leading zeros and string typing
$targetCells = @('C2','C4','C5','C6','C7','C9','C10','C11','C13','C14','C16','C18','C20','A23','B23','C23','A24','B24','C24','A25','B25','C25','A26','B26','C26','C29','C30','A33','B33','C33','A34','B34','C34','A35','B35','C35','A36','B36','C36','A37','B37','C37','A38','B38','C38','C39','C42','C43','A46','B46','C46','A47','B47','A48','B48','C48','A49','B49','A50','B50','C50','A51','B51','A52','B52','A53','B53','A54','B54','C54','A55','B55','C55','A56','B56','C56','A57','B57','A58','B58','C59','A60','B60','A61','B61','A62','B62','A63','B63','A64','B64','A66','B66','A67','B67','A68','B68','C68','A69','B69','A70','B70','A71','B71','A72','B72','A73','B73','A74','B74','A75','B75','A76','B76','A77','B77','A78','B78','A79','B79')
foreach ($cellref in $targetCells) {
    $ws.Range($cellref).NumberFormat = "@"
}

# Apply updated values
$ws.Range('C2').Value = '8603'
$ws.Range('C4').Value = '809'
$ws.Range('C5').Value = '604'
$ws.Range('C6').Value = '508'
$ws.Range('C7').Value = '443'
$ws.Range('C9').Value = '418'
$ws.Range('C10').Value = '376'
$ws.Range('C11').Value = '357'
$ws.Range('C13').Value = '249'
$ws.Range('C14').Value = '246'
$ws.Range('C16').Value = '200'
$ws.Range('C18').Value = '177'
$ws.Range('C20').Value = '167'
$ws.Range('A23').Value = 'Siscom SPA'
$ws.Range('B23').Value = '01778000040'
$ws.Range('C23').Value = '115'
$ws.Range('A24').Value = 'Regione Toscana'
$ws.Range('B24').Value = '01386030488'
$ws.Range('C24').Value = '113'
$ws.Range('A25').Value = 'Next Step Solution'
$ws.Range('B25').Value = '02554480349'
$ws.Range('C25').Value = '111'
$ws.Range('A26').Value = 'Regione Basilicata'
$ws.Range('B26').Value = '80002950766'
$ws.Range('C26').Value = '106'
$ws.Range('C29').Value = '80'
$ws.Range('C30').Value = '73'
$ws.Range('A33').Value = 'Regione Umbria'
$ws.Range('B33').Value = '80000130544'
$ws.Range('C33').Value = '57'
$ws.Range('A34').Value = 'NORDCOM'
$ws.Range('B34').Value = '13384100155'
$ws.Range('C34').Value = '56'
$ws.Range('A35').Value = 'DCS SOFTWARE E SERVIZI S.R.L.'
$ws.Range('B35').Value = '08063140019'
$ws.Range('C35').Value = '52'
$ws.Range('A36').Value = 'Italriscossioni Società Italiana di Fiscalità Locale S.r.l.'
$ws.Range('B36').Value = '06092371001'
$ws.Range('C36').Value = '50'
$ws.Range('A37').Value = 'Bluenext S.r.l.'
$ws.Range('B37').Value = '04228480408'
$ws.Range('C37').Value = '46'
$ws.Range('A38').Value = 'CINECA consorzio universitario'
$ws.Range('B38').Value = '00317740371'
$ws.Range('C38').Value = '44'
$ws.Range('C39').Value = '43'
$ws.Range('C42').Value = '35'
$ws.Range('C43').Value = '34'
$ws.Range('A46').Value = 'Comune di Palermo'
$ws.Range('B46').Value = '80016350821'
$ws.Range('C46').Value = '24'
$ws.Range('A47').Value = 'Citta'' Metropolitana di Roma Capitale'
$ws.Range('B47').Value = '80034390585'
$ws.Range('A48').Value = 'Nexi SpA'
$ws.Range('B48').Value = '13212880150'
$ws.Range('C48').Value = '19'
$ws.Range('A49').Value = 'ANDREANI TRIBUTI srl'
$ws.Range('B49').Value = '01412920439'
$ws.Range('A50').Value = 'Regione Lazio'
$ws.Range('B50').Value = '80143490581'
$ws.Range('C50').Value = '18'
$ws.Range('A51').Value = 'Servizi Locali SpA'
$ws.Range('B51').Value = '03170580751'
$ws.Range('A52').Value = 'Si.Form Consulting srl'
$ws.Range('B52').Value = '03943960827'
$ws.Range('A53').Value = 'Crédit Agricole Group Solutions Società Consortile per azioni'
$ws.Range('B53').Value = '02771790348'
$ws.Range('A54').Value = 'Aric Agenzia Regionale di Informatica e Committenza'
$ws.Range('B54').Value = '91022630676'
$ws.Range('C54').Value = '13'
$ws.Range('A55').Value = 'Comune di Catania'
$ws.Range('B55').Value = '00137020871'
$ws.Range('C55').Value = '9'
$ws.Range('A56').Value = 'UBI Banca'
$ws.Range('B56').Value = '03053920165'
$ws.Range('C56').Value = '9'
$ws.Range('A57').Value = 'Be Smart s.r.l.'
$ws.Range('B57').Value = '05817461006'
$ws.Range('A58').Value = 'ARCA Servizi s.r.l'
$ws.Range('B58').Value = '09106071005'
$ws.Range('C59').Value = '7'
$ws.Range('A60').Value = 'Phoenix IT Solutions S.r.L'
$ws.Range('B60').Value = '07623321218'
$ws.Range('A61').Value = 'CityPoste Payment Digital S.r.l.'
$ws.Range('B61').Value = '02003750672'
$ws.Range('A62').Value = 'ARGO SOFTWARE SRL'
$ws.Range('B62').Value = '00838520880'
$ws.Range('A63').Value = 'Linea Comune Spa'
$ws.Range('B63').Value = '05591950489'
$ws.Range('A64').Value = 'e-SED Società Cooperativa'
$ws.Range('B64').Value = '02695640421'
$ws.Range('A66').Value = 'Softline srl'
$ws.Range('B66').Value = '12299030150'
$ws.Range('A67').Value = 'KOINE'' SRL'
$ws.Range('B67').Value = '01934790971'
$ws.Range('A68').Value = 'ICCREA Banca SpA'
$ws.Range('B68').Value = '04774801007'
$ws.Range('C68').Value = '2'
$ws.Range('A69').Value = 'BANCA MONTE DEI PASCHI DI SIENA'
$ws.Range('B69').Value = '00884060526'
$ws.Range('A70').Value = 'Società Almaviva S.p.A.'
$ws.Range('B70').Value = '08450891000'
$ws.Range('A71').Value = 'Banca Nazionale del Lavoro S.p.A.'
$ws.Range('B71').Value = '09339391006'
$ws.Range('A72').Value = 'Banco BPM Società per Azioni'
$ws.Range('B72').Value = '09722490969'
$ws.Range('A73').Value = 'Engineering Ingegneria Informatica SpA'
$ws.Range('B73').Value = '00967720285'
$ws.Range('A74').Value = 'Noviservice srl'
$ws.Range('B74').Value = '02789990922'
$ws.Range('A75').Value = 'MegASP S.r.l.'
$ws.Range('B75').Value = '09898030151'
$ws.Range('A76').Value = 'I.C.A. - Imposte Comunali Affini – s.r.l.'
$ws.Range('B76').Value = '02478610583'
$ws.Range('A77').Value = 'Agenzia Italiana del Farmaco - AIFA'
$ws.Range('B77').Value = '97345810580'
$ws.Range('A78').Value = 'Ministero dello Sviluppo Economico'
$ws.Range('B78').Value = '80230390587'
$ws.Range('A79').Value = 'San Marco SPA'
$ws.Range('B79').Value = '04142440728'
